# DSM Scheduled Flights vs actual.xlsx
# Append 24 days of data (2022-05-03 .. 2022-05-26) to the Ark1 sheet,
# continuing the existing "Scheduled flights" / "Tracked flights" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$dates = @(
    "2022-05-03","2022-05-04","2022-05-05","2022-05-06","2022-05-07",
    "2022-05-08","2022-05-09","2022-05-10","2022-05-11","2022-05-12",
    "2022-05-13","2022-05-14","2022-05-15","2022-05-16","2022-05-17",
    "2022-05-18","2022-05-19","2022-05-20","2022-05-21","2022-05-22",
    "2022-05-23","2022-05-24","2022-05-25","2022-05-26"
)
$scheduled = @(58,71,79,81,47,62,61,57,57,82,73,45,58,61,66,72,71,68,56,57,62,59,66,81)
$tracked   = @(57,71,73,69,46,59,60,55,56,81,73,43,57,58,64,69,67,64,51,55,57,55,63,78)

$firstRow = 758
$lastRow  = $firstRow + $dates.Length - 1   # 781

# Carry the row-757 formatting (number formats, fonts, alignment) down onto
# the new rows before we populate them, so the new cells pick up the same
# styles already used throughout the table (text dates in col A, 0.00 in
# col B/C, 0.0% formula in col D).
$ws.Range("A757:D757").Copy()
$ws.Range("A$firstRow" + ":D$lastRow").PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $firstRow + $i
    $ws.Range("A$r").Value = $dates[$i]
    $ws.Range("B$r").Value = $scheduled[$i]
    $ws.Range("C$r").Value = $tracked[$i]
}

# Fill the ratio formula down the new rows in the same two chunks the
# workbook's author used (758:769, then 770:781), so each chunk becomes its
# own shared-formula group like the rest of column D.
$ws.Range("D$firstRow" + ":D769").Formula = "=C$firstRow/B$firstRow"
$ws.Range("D770:D$lastRow").Formula = "=C770/B770"

# Match the saved selection/scroll state left behind by the edit.
$ws.Range("A758:XFD758").Select() | Out-Null
